$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-28"

# Update the row label for November
$ws.Range("A12").Value = "November (through 11-28)"

# Update November row (row 12) values for columns C..I (B stays the same)
$ws.Range("C12").Value = 73
$ws.Range("D12").Value = 101
$ws.Range("E12").Value = 63
$ws.Range("F12").Value = 48
$ws.Range("G12").Value = 199
$ws.Range("H12").Value = 189
$ws.Range("I12").Value = 106

# Update Total row (row 13) values for columns C..I (B stays the same)
$ws.Range("C13").Value = 559
$ws.Range("D13").Value = 811
$ws.Range("E13").Value = 678
$ws.Range("F13").Value = 530
$ws.Range("G13").Value = 1256
$ws.Range("H13").Value = 1630
$ws.Range("I13").Value = 1504
